$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains text formatting so numeric-looking
# strings (e.g. "0.9991", "19.00") are not coerced into floating point numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.223.25"
$ws.Range("E2").Value = "  -0.53%  "

$ws.Range("D3").Value = "1.881.33"
$ws.Range("E3").Value = "  -1.38%  "

$ws.Range("D4").Value = "0.9991"
$ws.Range("E4").Value = "  -0.17%  "

$ws.Range("D5").Value = "237.69"
$ws.Range("E5").Value = "  -0.37%  "

$ws.Range("D6").Value = "0.9995"
$ws.Range("E6").Value = "  -0.20%  "

$ws.Range("D7").Value = "0.4668"
$ws.Range("E7").Value = "  -1.62%  "

$ws.Range("D8").Value = "0.2805"
$ws.Range("E8").Value = "  -1.31%  "

$ws.Range("D9").Value = "0.06543"
$ws.Range("E9").Value = "  -1.54%  "

$ws.Range("D10").Value = "19.49"
$ws.Range("E10").Value = "  +3.84%  "

$ws.Range("B11").Value = "Litecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D11").Value = "98.36"
$ws.Range("E11").Value = "  -2.82%  "

$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "0.07729"
$ws.Range("E12").Value = "  +0.15%  "

$ws.Range("D13").Value = "1.879.40"
$ws.Range("E13").Value = "  -1.49%  "

$ws.Range("D14").Value = "5.106"
$ws.Range("E14").Value = "  -1.46%  "

$ws.Range("D15").Value = "0.6657"
$ws.Range("E15").Value = "  -0.85%  "

$ws.Range("D16").Value = "286.09"
$ws.Range("E16").Value = "  +13.85%  "

$ws.Range("D17").Value = "30.208.16"
$ws.Range("E17").Value = "  -0.65%  "

$ws.Range("E18").Value = "  -0.16%  "

$ws.Range("D19").Value = "2.132.40"
$ws.Range("E19").Value = "  -1.14%  "

$ws.Range("D20").Value = "12.51"
$ws.Range("E20").Value = "  -1.26%  "

$ws.Range("D21").Value = "0.000007275"
$ws.Range("E21").Value = "  -2.49%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "5.299"
$ws.Range("E22").Value = "  -1.49%  "

$ws.Range("B23").Value = "BinanceUSD"
$ws.Range("C23").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D23").Value = "0.9985"
$ws.Range("E23").Value = "  -0.22%  "

$ws.Range("D24").Value = "6.185"
$ws.Range("E24").Value = "  -1.70%  "

$ws.Range("D25").Value = "166.59"
$ws.Range("E25").Value = "  +1.36%  "

$ws.Range("D26").Value = "9.226"
$ws.Range("E26").Value = "  -1.82%  "

$ws.Range("D27").Value = "19.00"
$ws.Range("E27").Value = "  +0.76%  "

$ws.Range("D28").Value = "1.983"
$ws.Range("E28").Value = "  -3.40%  "

$ws.Range("D29").Value = "1.371"
$ws.Range("E29").Value = "  -0.83%  "

$ws.Range("D30").Value = "0.09813"
$ws.Range("E30").Value = "  -2.75%  "

$ws.Range("D31").Value = "4.449"
$ws.Range("E31").Value = "  -3.78%  "

$ws.Range("D32").Value = "1.491"
$ws.Range("E32").Value = "  -1.38%  "

$ws.Range("D33").Value = "4.173"
$ws.Range("E33").Value = "  -1.47%  "

$ws.Range("D34").Value = "0.04668"
$ws.Range("E34").Value = "  -1.51%  "

$ws.Range("D35").Value = "0.7073"
$ws.Range("E35").Value = "  -2.56%  "

$ws.Range("D36").Value = "1.092"
$ws.Range("E36").Value = "  -1.48%  "

$ws.Range("D37").Value = "0.9992"
$ws.Range("E37").Value = "  -0.15%  "

$ws.Range("D38").Value = "2.710"
$ws.Range("E38").Value = "  -0.41%  "

$ws.Range("D39").Value = "0.01862"
$ws.Range("E39").Value = "  -2.74%  "

$ws.Range("D40").Value = "6.715"
$ws.Range("E40").Value = "  +7.92%  "

$ws.Range("D41").Value = "2.504"
$ws.Range("E41").Value = "  -3.69%  "

$ws.Range("D42").Value = "72.35"
$ws.Range("E42").Value = "  -1.27%  "

$ws.Range("D43").Value = "0.8680"
$ws.Range("E43").Value = "  +0.76%  "

$ws.Range("D44").Value = "1.950"
$ws.Range("E44").Value = "  -1.47%  "

$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "104.06"
$ws.Range("E45").Value = "  -2.56%  "

$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").Value = "0.9992"
$ws.Range("E46").Value = "  -0.13%  "

$ws.Range("D47").Value = "0.4165"
$ws.Range("E47").Value = "  -1.47%  "

$ws.Range("D48").Value = "993.48"
$ws.Range("E48").Value = "  -5.83%  "

$ws.Range("D49").Value = "7.224"
$ws.Range("E49").Value = "  -3.04%  "

$ws.Range("D50").Value = "9.318"
$ws.Range("E50").Value = "  +5.89%  "

$ws.Range("D51").Value = "0.1156"
$ws.Range("E51").Value = "  -3.23%  "

